$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 22-25: re-key the existing rows onto the "copied" style (index 2)
# that the rest of the table already uses. The original file had these four
# rows authored with a slightly different (but visually identical) General
# style than their neighbours; re-typing them picks up the same style as the
# row immediately above once the blank row is (re)inserted there.
$ws.Range("A22:J25").Delete()

$ws.Rows.Item(22).Insert()
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(25).Insert()

$ws.Range("A22").Value = "24 Apr 2020"
$ws.Range("B22").Value = "0,6"
$ws.Range("C22").Value = "-7.5,6"
$ws.Range("D22").Value = "-20,6"
$ws.Range("E22").Value = "-32.5,6"
$ws.Range("F22").Value = "-25,6"
$ws.Range("G22").Value = "-15,6"
$ws.Range("H22").Value = "-27.5,6"
$ws.Range("I22").Value = "-12.5,6"
$ws.Range("J22").Value = "-22.5,6"

$ws.Range("A23").Value = "28 Apr 2020"
$ws.Range("B23").Value = "0,6"
$ws.Range("C23").Value = "-7.5,6"
$ws.Range("D23").Value = "-7.5,6"
$ws.Range("E23").Value = "-30,6"
$ws.Range("F23").Value = "-25,6"
$ws.Range("G23").Value = "0,6"
$ws.Range("H23").Value = "-25,6"
$ws.Range("I23").Value = "-10,6"
$ws.Range("J23").Value = "-22.5,6"

$ws.Range("A24").Value = "1 May 2020"
$ws.Range("B24").Value = "0,6"
$ws.Range("C24").Value = "-5,6"
$ws.Range("D24").Value = "-7.5,6"
$ws.Range("E24").Value = "-30,6"
$ws.Range("F24").Value = "-25,6"
$ws.Range("G24").Value = "0,6"
$ws.Range("H24").Value = "-25,6"
$ws.Range("I24").Value = "-10,6"
$ws.Range("J24").Value = "-22.5,6"

$ws.Range("A25").Value = "7 May 2020"
$ws.Range("B25").Value = "0,6"
$ws.Range("C25").Value = "-5,6"
$ws.Range("D25").Value = "-7.5,6"
$ws.Range("E25").Value = "-30,6"
$ws.Range("F25").Value = "-25,6"
$ws.Range("G25").Value = "0,6"
$ws.Range("H25").Value = "-25,6"
$ws.Range("I25").Value = "-10,6"
$ws.Range("J25").Value = "-22.5,6"

# --- New row 26 (18 May 2020): inherit formatting from row above (25, all
# style 2), then reset the date/jug columns back to the plain "typed in"
# look (no explicit style) to match how the author entered them.
$ws.Rows.Item(26).Insert()
$ws.Range("A26:B26").ClearFormats()

$ws.Range("A26").Value = "18 May 2020"
$ws.Range("B26").Value = "0,6"
$ws.Range("C26").Value = "-5,6"
$ws.Range("D26").Value = "-5,6"
$ws.Range("E26").Value = "-27.5,6"
$ws.Range("F26").Value = "-22.5,6"
$ws.Range("G26").Value = "0,6"
$ws.Range("H26").Value = "-22.5,6"
$ws.Range("I26").Value = "-5,6"
$ws.Range("J26").Value = "-20,6"

# --- New row 27 (2 Jun 2020)
$ws.Rows.Item(27).Insert()
$ws.Range("A27:B27").ClearFormats()

$ws.Range("A27").Value = "2 Jun 2020"
$ws.Range("B27").Value = "0,6"
$ws.Range("C27").Value = "-5,6"
$ws.Range("D27").Value = "-5,6"
$ws.Range("E27").Value = "-22.5,6"
$ws.Range("F27").Value = "-22.5,6"
$ws.Range("G27").Value = "0,6"
$ws.Range("H27").Value = "-20,6"
$ws.Range("I27").Value = "0,6"
$ws.Range("J27").Value = "-20,6"

$ws.Range("D31").Select() | Out-Null
